$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("May 08")

# Update "Last Updated" timestamp
$ws.Range("B1").Value = "May 08 2022 22:35"

# Update Balance and Total Consumed values
$ws.Range("B3").Value = -1611
$ws.Range("B4").Value = 1289

# New food-log rows use text values (e.g. "20.00") that look numeric, so
# pre-format the range as Text to stop Excel auto-converting them to
# numbers, then restore the default "Normal" style afterwards so the
# saved cells carry no explicit style (matching the rest of the sheet).
$newRange = $ws.Range("B18:G19")
$newRange.NumberFormat = "@"

$ws.Range("A18").Value = "פירות קפואים: תות שדה"
$ws.Range("B18").Value = "20.00"
$ws.Range("C18").Value = "גרם"
$ws.Range("D18").Value = "6.60"
$ws.Range("E18").Value = "0.16"
$ws.Range("F18").Value = "1.14"
$ws.Range("G18").Value = "0.00"

$ws.Range("A19").Value = "בולגרית פיראוס 5%"
$ws.Range("B19").Value = "80.00"
$ws.Range("C19").Value = "גרם"
$ws.Range("D19").Value = "95.20"
$ws.Range("E19").Value = "14.40"
$ws.Range("F19").Value = "0.40"
$ws.Range("G19").Value = "4.00"

$newRange.Style = "Normal"
